$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matches")

# home-team / away-team "evaluated" lookup formulas for each group-stage
# match row (2-49), pulling the team name from the seeds table via the
# home-seed / away-seed columns.
$homeFormula = "=INDEX(seeds[team],MATCH(matches[[#This Row],[home-seed]],seeds[seed],0))"
$awayFormula = "=INDEX(seeds[team],MATCH(matches[[#This Row],[away-seed]],seeds[seed],0))"

for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 6).Formula = $homeFormula
    $ws.Cells.Item($r, 7).Formula = $awayFormula
}

# Matches becomes the active/selected sheet (Tournament loses tabSelected).
$ws.Activate()
$ws.Range("F50:G65").Select()
